$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-25, columns G,H,I,J,M,N,O,P,Q,R,S,T
$data = @{}
$data[2] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"22.32577533333334"; 14=[double]"66.977326"; 15=[double]"0.1836188216937888"; 16=[double]"0.1836188216937889"; 17=[double]"523.4377153342481"; 18=[double]"4710.939438008232"; 19=[double]"0.1631274145983724"; 20=[double]"0.1631274145983725"}
$data[3] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"52.60848633333333"; 14=[double]"157.825459"; 15=[double]"0.4326796327291624"; 16=[double]"0.4326796327291624"; 17=[double]"1233.429320252932"; 18=[double]"11100.86388227639"; 19=[double]"0.3843936541221641"; 20=[double]"0.3843936541221641"}
$data[4] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"0.127903"; 14=[double]"0.383709"; 15=[double]"0.001051940987511236"; 16=[double]"0.001051940987511237"; 17=[double]"2.998742623931999"; 18=[double]"26.988683615388"; 19=[double]"0.0009345469708379649"; 20=[double]"0.0009345469708379652"}
$data[5] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"46.272583"; 14=[double]"138.817749"; 15=[double]"0.3805698588439337"; 16=[double]"0.3805698588439337"; 17=[double]"1084.881253461852"; 18=[double]"9763.931281156667"; 19=[double]"0.3380992023290956"; 20=[double]"0.3380992023290957"}
$data[6] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"0.1965033333333333"; 14=[double]"0.58951"; 15=[double]"0.001616145911479139"; 16=[double]"0.001616145911479139"; 17=[double]"4.60710789748"; 18=[double]"41.46397107732"; 19=[double]"0.001435788018468915"; 20=[double]"0.001435788018468915"}
$data[7] = @{7=[double]"23.445444"; 8=[double]"70.336332"; 9=[double]"0.8884024692763315"; 10=[double]"0.8884024692763316"; 13=[double]"0.056368"; 14=[double]"0.169104"; 15=[double]"0.0004635998341245583"; 16=[double]"0.0004635998341245584"; 17=[double]"1.321572787392"; 18=[double]"11.894155086528"; 19=[double]"0.0004118632373923553"; 20=[double]"0.0004118632373923554"}
$data[8] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"22.32577533333334"; 14=[double]"66.977326"; 15=[double]"0.1836188216937888"; 16=[double]"0.1836188216937889"; 17=[double]"55.82828031404001"; 18=[double]"502.4545228263601"; 19=[double]"0.01739867564431631"; 20=[double]"0.01739867564431632"}
$data[9] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"52.60848633333333"; 14=[double]"157.825459"; 15=[double]"0.4326796327291624"; 16=[double]"0.4326796327291624"; 17=[double]"131.55383309486"; 18=[double]"1183.98449785374"; 19=[double]"0.04099826215152786"; 20=[double]"0.04099826215152787"}
$data[10] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"0.127903"; 14=[double]"0.383709"; 15=[double]"0.001051940987511236"; 16=[double]"0.001051940987511237"; 17=[double]"0.31983679986"; 18=[double]"2.87853119874"; 19=[double]"9.96759475408882e-05"; 20=[double]"9.967594754088824e-05"}
$data[11] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"46.272583"; 14=[double]"138.817749"; 15=[double]"0.3805698588439337"; 16=[double]"0.3805698588439337"; 17=[double]"115.71014650146"; 18=[double]"1041.39131851314"; 19=[double]"0.03606063622971623"; 20=[double]"0.03606063622971625"}
$data[12] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"0.1965033333333333"; 14=[double]"0.58951"; 15=[double]"0.001616145911479139"; 16=[double]"0.001616145911479139"; 17=[double]"0.4913801654"; 18=[double]"4.4224214886"; 19=[double]"0.0001531367985500184"; 20=[double]"0.0001531367985500184"}
$data[13] = @{7=[double]"2.50062"; 8=[double]"7.501860000000001"; 9=[double]"0.0947543148563013"; 10=[double]"0.09475431485630131"; 13=[double]"0.056368"; 14=[double]"0.169104"; 15=[double]"0.0004635998341245583"; 16=[double]"0.0004635998341245584"; 17=[double]"0.14095494816"; 18=[double]"1.26859453344"; 19=[double]"4.392808464996745e-05"; 20=[double]"4.392808464996747e-05"}
$data[14] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"22.32577533333334"; 14=[double]"66.977326"; 15=[double]"0.1836188216937888"; 16=[double]"0.1836188216937889"; 17=[double]"6.279518734530889"; 18=[double]"56.515668610778"; 19=[double]"0.001956988627447198"; 20=[double]"0.001956988627447199"}
$data[15] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"52.60848633333333"; 14=[double]"157.825459"; 15=[double]"0.4326796327291624"; 16=[double]"0.4326796327291624"; 17=[double]"14.79706619784189"; 18=[double]"133.173595780577"; 19=[double]"0.004611450573357228"; 20=[double]"0.00461145057335723"}
$data[16] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"0.127903"; 14=[double]"0.383709"; 15=[double]"0.001051940987511236"; 16=[double]"0.001051940987511237"; 17=[double]"0.03597497836966666"; 18=[double]"0.323774805327"; 19=[double]"1.12114680309742e-05"; 20=[double]"1.12114680309742e-05"}
$data[17] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"46.272583"; 14=[double]"138.817749"; 15=[double]"0.3805698588439337"; 16=[double]"0.3805698588439337"; 17=[double]"13.01498145104966"; 18=[double]"117.134833059447"; 19=[double]"0.004056070498855383"; 20=[double]"0.004056070498855384"}
$data[18] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"0.1965033333333333"; 14=[double]"0.58951"; 15=[double]"0.001616145911479139"; 16=[double]"0.001616145911479139"; 17=[double]"0.05527003405888888"; 18=[double]"0.49743030653"; 19=[double]"1.72247002779179e-05"; 20=[double]"1.72247002779179e-05"}
$data[19] = @{7=[double]"0.2812676666666666"; 8=[double]"0.843803"; 9=[double]"0.0106578868625503"; 10=[double]"0.0106578868625503"; 13=[double]"0.056368"; 14=[double]"0.169104"; 15=[double]"0.0004635998341245583"; 16=[double]"0.0004635998341245584"; 17=[double]"0.01585449583466667"; 18=[double]"0.142690462512"; 19=[double]"4.940994581596628e-06"; 20=[double]"4.94099458159663e-06"}
$data[20] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"22.32577533333334"; 14=[double]"66.977326"; 15=[double]"0.1836188216937888"; 16=[double]"0.1836188216937889"; 17=[double]"3.644333052686445"; 18=[double]"32.798997474178"; 19=[double]"0.001135742823652885"; 20=[double]"0.001135742823652885"}
$data[21] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"52.60848633333333"; 14=[double]"157.825459"; 15=[double]"0.4326796327291624"; 16=[double]"0.4326796327291624"; 17=[double]"8.587511194297445"; 18=[double]"77.287600748677"; 19=[double]"0.002676265882113189"; 20=[double]"0.00267626588211319"}
$data[22] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"0.127903"; 14=[double]"0.383709"; 15=[double]"0.001051940987511236"; 16=[double]"0.001051940987511237"; 17=[double]"0.02087816093633333"; 18=[double]"0.187903448427"; 19=[double]"6.506601101408928e-06"; 20=[double]"6.506601101408929e-06"}
$data[23] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"46.272583"; 14=[double]"138.817749"; 15=[double]"0.3805698588439337"; 16=[double]"0.3805698588439337"; 17=[double]"7.553274237616334"; 18=[double]"67.97946813854699"; 19=[double]"0.002353949786266437"; 20=[double]"0.002353949786266437"}
$data[24] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"0.1965033333333333"; 14=[double]"0.58951"; 15=[double]"0.001616145911479139"; 16=[double]"0.001616145911479139"; 17=[double]"0.03207609061444445"; 18=[double]"0.28868481553"; 19=[double]"9.996394182288082e-06"; 20=[double]"9.996394182288084e-06"}
$data[25] = @{7=[double]"0.1632343333333333"; 8=[double]"0.489703"; 9=[double]"0.006185329004816848"; 10=[double]"0.006185329004816848"; 13=[double]"0.056368"; 14=[double]"0.169104"; 15=[double]"0.0004635998341245583"; 16=[double]"0.0004635998341245584"; 17=[double]"0.009201192901333335"; 18=[double]"0.082810736112"; 19=[double]"2.86751750063891e-06"; 20=[double]"2.86751750063891e-06"}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item($r, $c).Value2 = $data[$r][$c]
    }
}
